$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# G column literal values for the new rows (the fill-down series keeps using
# epsilon/10 rounded to 3dp, same pattern already present from row 13 on:
# the formula chain for G stopped being "shared" at row 13, continuing as a
# plain literal of 0.001/0.002 thereafter).
$gValues = @{
    14 = 0.001
    15 = 0.002
    16 = 0.002
    17 = 0.002
    18 = 0.002
    19 = 0.002
    20 = 0.002
    21 = 0.002
}

for ($r = 14; $r -le 21; $r++) {
    $p = $r - 1
    $ws.Range("A$r").Formula = "=A$p+1"
    $ws.Range("B$r").Formula = "=B$p-(B$p*G$p)"
    $ws.Range("C$r").Formula = "=IF(H$r=TRUE,B$r,IF(D$r>epsilon,B$r,C$p))"
    $ws.Range("D$r").Formula = "=(C$p-B$r)/C$p"
    $ws.Range("E$r").Formula = "=IF(D$r>=epsilon,TRUE,FALSE)"
    $ws.Range("F$r").Formula = "=IF(E$r=TRUE,0,MOD(F$p+1,5))"
    $ws.Range("G$r").Value = $gValues[$r]
    $ws.Range("H$r").Formula = "=IF(F$p=patience-1,TRUE)"
    $ws.Range("I$r").Formula = "=IF(C$r<>C$p,A$r,I$p)"
}

[void]$ws.Range("H4:H21").Select()
